$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 into I1 and J1, then set header labels
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I and J column values for data rows 2-61
$iValues = @(6,7,9,10,5,8,5,7,9,5,6,9,5,8,8,5,6,8,6,1,1,1,1,7,6,5,6,6,8,8,5,6,8,9,3,8,4,7,4,6,1,4,1,6,8,7,8,8,6,6,5,9,9,8,9,7,6,5,6,4)
$jValues = @(6,8,9,10,5,8,5,7,9,5,6,9,5,8,8,5,6,8,6,2,1,3,3,7,6,5,7,6,8,8,5,6,8,9,3,8,5,7,5,6,2,4,1,6,8,7,8,8,6,6,6,9,9,8,9,7,6,5,6,4)

For ($k = 0; $k -lt 60; $k++) {
    $r = $k + 2
    $ws.Cells.Item($r, 9).Value = $iValues[$k]
    $ws.Cells.Item($r, 10).Value = $jValues[$k]
}

$ws.Application.CutCopyMode = $false
